# This workbook holds a weekly price report. This edit refreshes the data
# rows (2-13) by re-ordering which week's figures sit in which row: for
# every destination row we copy the full record (Fecha, Volumen, Precio
# minimo/maximo/promedio ponderado, Origen, Precio $/Kg) that used to sit
# in a different source row. The other columns (A, B, C, E, F, G, H, I, N,
# Q, R) are identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to copy from, read from
# the original/"before" state of the sheet).
$rowMap = @{
    2  = 6
    3  = 7
    4  = 12
    5  = 2
    6  = 13
    7  = 5
    8  = 3
    9  = 8
    10 = 11
    11 = 10
    12 = 9
    13 = 4
}

# Columns whose values move together with the row record.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot all source values first, since rows 2-13 both provide and
# receive data (the mapping is a permutation), so writes must not clobber
# values still needed as a source for another row.
$snapshot = @{}
foreach ($r in 2..13) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
